$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "x" mark from C13 to D13
$ws.Range("C13").Value = $null
$ws.Range("D13").Value = "x"

# Update the active selection to C13
$ws.Range("C13").Select()
